# Rename the inline picture shapes that live in the document's headers
# and footers:
#   - the two Pearson/Edexcel logo pictures (descr contains
#     "PearsonLogo.png"), currently named "image2.png", are renamed to
#     "image1.png"
#   - the BTEC logo picture (descr "BTec_Logo-Orange"), currently named
#     "image1.jpg", is renamed to "image2.jpg"
#
# Word's InlineShape object has no writable "Name" in the real object
# model either - renaming is done by flipping the picture to a floating
# Shape (ConvertToShape), setting its Name, then converting it back to
# an inline shape (ConvertToInlineShape). That is the route used below.

$d = $word.ActiveDocument

function Rename-InlineLogo($inlineShape, $oldName, $newName) {
    if ($inlineShape.Width -eq $null) { return }
    $floatingShape = $inlineShape.ConvertToShape()
    if ($floatingShape.Name -eq $oldName) {
        $floatingShape.Name = $newName
    }
    $floatingShape.ConvertToInlineShape() | Out-Null
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers.Item($h)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shape = $shapes.Item($i)
                if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-InlineLogo $shape "image1.jpg" "image2.jpg"
                }
            }
        }
    }

    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shape = $shapes.Item($i)
                if ($shape.AlternativeText -like "*PearsonLogo.png") {
                    Rename-InlineLogo $shape "image2.png" "image1.png"
                }
            }
        }
    }
}
